$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.556.81'
$ws.Range('E2').Value = '  +2.77%  '

$ws.Range('D3').Value = '2.193.74'
$ws.Range('E3').Value = '  +1.63%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.94%  '

$ws.Range('E6').Value = '  +1.16%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.22'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.34%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('E9').Value = '  +2.47%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.05'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.49%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0919'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.12%  '

$ws.Range('E12').Value = '  +2.27%  '

$ws.Range('E13').Value = '  +1.48%  '

$ws.Range('D14').Value = '2.524.27'
$ws.Range('E14').Value = '  +1.69%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.40%  '

$ws.Range('D16').Value = '2.199.41'
$ws.Range('E16').Value = '  +2.41%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.776'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.09%  '

$ws.Range('D18').Value = '42.435.77'
$ws.Range('E18').Value = '  +2.64%  '

$ws.Range('E19').Value = '  +1.11%  '

$ws.Range('E20').Value = '  +2.25%  '

$ws.Range('E21').Value = '  +2.60%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.57%  '

$ws.Range('E23').Value = '  +7.91%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.40'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.13%  '

$ws.Range('E25').Value = '  -0.15%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.65'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.94%  '

$ws.Range('E28').Value = '  +2.61%  '

$ws.Range('E29').Value = '  +1.11%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +13.55%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.55'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.84%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0796'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.04%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.23%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.120'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.49%  '

$ws.Range('E36').Value = '  +3.63%  '

$ws.Range('E37').Value = '  +3.80%  '

$ws.Range('E38').Value = '  +10.00%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.11'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.97%  '

$ws.Range('E40').Value = '  +0.56%  '

$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.25'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.87%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.197'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.85%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.97'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.00%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.483'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +25.55%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.47%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.41'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.64%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0976'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.46%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.42'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +12.52%  '

$ws.Range('E49').Value = '  +2.47%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.46%  '

$ws.Range('E51').Value = '  +1.53%  '
